$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Insulin Receptors")

# --- Wipe the old contents of the sheet (rows 1-34, cols A-F) ---
$ws.Range("A1:F34").Clear()

# --- Title / body mass (unchanged from before) ---
$ws.Range("A1").Value = "Insulin Receptors (Units are mU, mU/mL and mU/Min)"
$ws.Range("A3").Value = "Body Mass"
$ws.Range("B3").Formula = "=BodyMassMale"

# --- Explanatory notes (rows 5-9) ---
$ws.Range("A5").Value = "Binding = K1 * [I] * Free - K2 * Bound"
$ws.Range("A6").Value = "[I] = Pool for General and Kidney and Portal for Liver"
$ws.Range("A7").Value = "K1 = 10 * K2"
$ws.Range("A8").Value = "General and Liver Degradation = K3 * Bound"
$ws.Range("A9").Value = "Kidney Degradation = K4 * [I] * GFR"

# blank, but numeric-formatted, cells left over beneath the notes
$ws.Range("B6").NumberFormat = "0"
$ws.Range("B7").NumberFormat = "0"
$ws.Range("B8").NumberFormat = "0"
$ws.Range("B9").NumberFormat = "0"
$ws.Range("B10").NumberFormat = "0"

# --- Receptors table (rows 11-13) ---
$ws.Range("C11").Value = "Total (mU)"
$ws.Range("D11").Value = "% Bound"
$ws.Range("E11").Value = "Bound"
$ws.Range("F11").Value = "Free"

$ws.Range("A12").Value = "Hepatic Receptors (mU)"
$ws.Range("C12").Value = 50000
$ws.Range("D12").Value = 20
$ws.Range("E12").Formula = "=0.01*D12*C12"
$ws.Range("F12").Formula = "=C12-E12"

$ws.Range("A13").Value = "General Receptors (mU)"
$ws.Range("C13").Value = 19000
$ws.Range("D13").Value = 8
$ws.Range("E13").Formula = "=0.01*D13*C13"
$ws.Range("F13").Formula = "=C13-E13"

# --- Secretion / ECFV / Portal vein values ---
$ws.Range("A15").Value = "Secretion=Degradation (mU/Min)"
$ws.Range("B15").Value = 17

$ws.Range("A17").Value = "[Insulin] ECFV"
$ws.Range("B17").Value = 0.02
$ws.Range("C17").Value = 15000
$ws.Range("D17").Formula = "=B17*C17"

$ws.Range("A18").Value = "[Insulin] Portal Vein"
$ws.Range("B18").Value = 0.052

# --- Degradation table (rows 20-24) ---
$ws.Range("A20").Value = "Degradation"
$ws.Range("B20").Value = "% Total"
$ws.Range("C20").Value = "mU/Min"
$ws.Range("D20").Value = "K3 and K4"
$ws.Range("E20").Value = "K1"
$ws.Range("F20").Value = "K2"

$ws.Range("A21").Value = "Liver"
$ws.Range("B21").Value = 79
$ws.Range("C21").Formula = "=0.01*B21*B15"
$ws.Range("D21").Formula = "=C21/E12"
$ws.Range("E21").Formula = "=10*F21"
$ws.Range("F21").Formula = "=C21/(10*B18*F12-E12)"

$ws.Range("A22").Value = "Kidney"
$ws.Range("B22").Value = 9
$ws.Range("C22").Formula = "=0.01*B22*B15"
$ws.Range("D22").Formula = "=C22/(1000*B17*125)"

$ws.Range("A23").Value = "General"
$ws.Range("B23").Value = 12
$ws.Range("C23").Formula = "=0.01*B23*B15"
$ws.Range("D23").Formula = "=C23/E13"
$ws.Range("E23").Formula = "=10*F23"
$ws.Range("F23").Formula = "=C23/(10*B17*F13-E13)"

$ws.Range("C24").Formula = "=SUM(C21:C23)"

# --- Column F needs to be shown/sized for the new data (best-fit width of 12) ---
$ws.Columns.Item(6).ColumnWidth = 11.166666666666666

# --- Sheet view: scroll back to the top and select D15 ---
$ws.Activate()
$ws.Range("D15").Select()

# --- Page setup: portrait, like the other sheets in the workbook ---
$ws.PageSetup.Orientation = 1
